$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the password value in B3 to "Test@123" (was "Test@1234")
$ws.Range("B3").Value = "Test@123"

# Move the active selection to E6 (was K11)
$ws.Range("E6").Select()
